$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple text/value updates (non-numeric-looking strings, safe to assign directly)
$ws.Range('D2').Value = '22.457.82'
$ws.Range('E2').Value = '  +0.28%  '
$ws.Range('D3').Value = '1.573.17'
$ws.Range('E3').Value = '  +0.05%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('E5').Value = '  -0.29%  '
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('E7').Value = '  -0.42%  '
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('E9').Value = '  -0.18%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('E10').Value = '  -1.00%  '
$ws.Range('B11').Value = 'Polygon'
$ws.Range('C11').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('E11').Value = '  -1.21%  '
$ws.Range('E12').Value = '  -0.11%  '
$ws.Range('E13').Value = '  +0.73%  '
$ws.Range('E14').Value = '  +0.53%  '
$ws.Range('E15').Value = '  +0.70%  '
$ws.Range('D16').Value = '1.572.74'
$ws.Range('E16').Value = '  -0.61%  '
$ws.Range('E17').Value = '  -0.84%  '
$ws.Range('E19').Value = '  -0.19%  '
$ws.Range('E20').Value = '  -0.23%  '
$ws.Range('E21').Value = '  +1.02%  '
$ws.Range('E22').Value = '  -1.79%  '
$ws.Range('E23').Value = '  +1.76%  '
$ws.Range('D24').Value = '22.453.35'
$ws.Range('E24').Value = '  +0.27%  '
$ws.Range('E25').Value = '  -4.86%  '
$ws.Range('E26').Value = '  -5.32%  '
$ws.Range('E27').Value = '  -0.36%  '
$ws.Range('E28').Value = '  +2.39%  '
$ws.Range('E29').Value = '  -1.12%  '
$ws.Range('E30').Value = '  +0.10%  '
$ws.Range('D31').Value = '1.747.73'
$ws.Range('E31').Value = '  -0.35%  '
$ws.Range('E32').Value = '  +3.60%  '
$ws.Range('E33').Value = '  -0.25%  '
$ws.Range('E34').Value = '  -1.80%  '
$ws.Range('E35').Value = '  -1.56%  '
$ws.Range('E36').Value = '  -1.22%  '
$ws.Range('E37').Value = '  +4.56%  '
$ws.Range('E38').Value = '  -3.27%  '
$ws.Range('E39').Value = '  -0.68%  '
$ws.Range('E40').Value = '  +0.36%  '
$ws.Range('E41').Value = '  +1.13%  '
$ws.Range('E42').Value = '  -1.15%  '
$ws.Range('E43').Value = '  -2.50%  '
$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('E44').Value = '  -0.17%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('E45').Value = '  -0.66%  '
$ws.Range('E46').Value = '  +0.71%  '
$ws.Range('E47').Value = '  -2.39%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('E48').Value = '  +3.90%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('E49').Value = '  +0.32%  '
$ws.Range('E50').Value = '  -4.55%  '
$ws.Range('E51').Value = '  -0.08%  '

# Price column updates whose text looks like a plain number; force text storage
# so Excel keeps them as strings (matching the source inlineStr cells) instead of
# silently converting to numeric values (which would also mangle trailing zeros).
$textCells = [ordered]@{
    'D4' = '1.002'
    'D6' = '291.11'
    'D8' = '49.91'
    'D9' = '0.3405'
    'D10' = '0.07578'
    'D11' = '1.146'
    'D13' = '21.39'
    'D14' = '5.992'
    'D15' = '6.957'
    'D18' = '91.15'
    'D19' = '0.06741'
    'D22' = '16.43'
    'D23' = '12.22'
    'D25' = '2.324'
    'D27' = '20.17'
    'D28' = '148.62'
    'D29' = '5.005'
    'D30' = '126.02'
    'D32' = '1.044'
    'D33' = '6.145'
    'D34' = '1.983'
    'D35' = '9.876'
    'D36' = '0.08451'
    'D37' = '1.386'
    'D38' = '0.02466'
    'D39' = '0.2300'
    'D40' = '0.06546'
    'D41' = '5.491'
    'D42' = '11.40'
    'D43' = '0.6285'
    'D44' = '1.001'
    'D45' = '14.01'
    'D46' = '3.815'
    'D47' = '0.5874'
    'D48' = '130.25'
    'D49' = '2.093'
    'D50' = '1.232'
    'D51' = '0.07328'
}
foreach ($ref in $textCells.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $textCells[$ref]
    $cell.ClearFormats()
}

